$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered/top-aligned) from
# the existing "IP" header cell (H1) onto the new header cells so they reuse
# the same cell style instead of creating a brand new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Add the two new data values in row 2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
